# Read multiple values from one column
#
# Adds a new column M ("Site") holding a hyperlinked value per row, mirroring
# the existing domain used in each row's e-mail address (column C):
#   row 1 -> patrick@maatwebsite.nl  => "Maatwebsite" -> https://www.maatwebsite.nl
#   row 2 -> taylor@laravel.com      => "Laravel"     -> https://laravel.com

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the extra column.
$ws.Range("M1").Value = "Maatwebsite"
$ws.Range("M2").Value = "Laravel"

# Turn them into real hyperlinks (adds the relationship + display text).
$ws.Hyperlinks.Add($ws.Range("M1"), "https://www.maatwebsite.nl")
$ws.Hyperlinks.Add($ws.Range("M2"), "https://laravel.com")

# Match the look of the other hyperlink cells (C1/C2) in the sheet.
$ws.Range("M1").Style = "Hyperlink"
$ws.Range("M2").Style = "Hyperlink"

# Size the new column to fit its contents, like the other autosized columns.
$ws.Columns.Item(13).ColumnWidth = 10

# Leave the selection where the author ended up after the edit.
$ws.Range("F10").Select() | Out-Null
